$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a text value into a cell while forcing text storage so that
# Excel does not auto-convert date/time-looking strings into date serials,
# then reset the cell style back to Normal so no stray style index remains.
function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Step 1: insert 2 new rows right after row 2 ---
# (old row 3 -> row 5, old row 4 -> row 6, old row 5 -> row 7, old row 6 -> row 8)
$ws.Rows("3:4").Insert()

# New row 3: Dinheiro, 03/05/2021, 22:00, 8.555, FALSE
Set-TextCell 3 1 "Dinheiro"
Set-TextCell 3 2 "03/05/2021"
Set-TextCell 3 3 "22:00"
$ws.Cells.Item(3, 4).Value = 8.555
$ws.Cells.Item(3, 5).Value = $false

# New row 4: Dinheiro, 03/05/2021, 21:25, 13.055, TRUE
Set-TextCell 4 1 "Dinheiro"
Set-TextCell 4 2 "03/05/2021"
Set-TextCell 4 3 "21:25"
$ws.Cells.Item(4, 4).Value = 13.055
$ws.Cells.Item(4, 5).Value = $true

# --- Step 2: insert 3 new rows right after (the now-shifted) row 6 ---
# (old row 5 [-> row 7] -> row 10, old row 6 [-> row 8] -> row 11)
$ws.Rows("7:9").Insert()

# New row 7: Cartão, 02/05/2021, 22:36, 8.555, TRUE
Set-TextCell 7 1 "Cartão"
Set-TextCell 7 2 "02/05/2021"
Set-TextCell 7 3 "22:36"
$ws.Cells.Item(7, 4).Value = 8.555
$ws.Cells.Item(7, 5).Value = $true

# New row 8: Cartão, 02/05/2021, 19:56, 60, TRUE
Set-TextCell 8 1 "Cartão"
Set-TextCell 8 2 "02/05/2021"
Set-TextCell 8 3 "19:56"
$ws.Cells.Item(8, 4).Value = 60
$ws.Cells.Item(8, 5).Value = $true

# New row 9: Cartão, 02/05/2021, 19:59, 160, TRUE
Set-TextCell 9 1 "Cartão"
Set-TextCell 9 2 "02/05/2021"
Set-TextCell 9 3 "19:59"
$ws.Cells.Item(9, 4).Value = 160
$ws.Cells.Item(9, 5).Value = $true
